$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("20-09-2021", 3631, 6605),
    @("21-09-2021", 3596, 6880),
    @("22-09-2021", 3804, 8512),
    @("23-09-2021", 3686, 9383),
    @("24-09-2021", 3293, 8356),
    @("27-09-2021", 3297, 7818),
    @("28-09-2021", 2835, 8996),
    @("29-09-2021", 2998, 8203),
    @("30-09-2021", 3740, 8903)
)

$startRow = 181
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
